# "Fruta / hortaliza, semanal" weekly price-sheet update.
#
# A new daily price observation is inserted as row 362 (Murcott / Primera,
# recorded 2021-11-05, i.e. Excel serial 44505), pushing the existing rows
# 362-427 down to 363-428 and growing the used range to A1:T428.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 362..427 down one, creating a blank row 362 with D362 inheriting
# the date-formatted style from the row that used to be there (matches how
# Excel's own Insert handles formatting continuity).
$ws.Rows.Item(362).Insert()

# Populate the new row with the inserted observation.
$ws.Range("A362").Value = 10
$ws.Range("B362").Value = "Vega Modelo de Temuco"
$ws.Range("C362").Value = "La Araucanía"
$ws.Range("D362").Value = 44505
$ws.Range("E362").Value = 9
$ws.Range("F362").Value = "Fruta"
$ws.Range("G362").Value = 100102
$ws.Range("H362").Value = "Cítricos"
$ws.Range("I362").Value = 100102004
$ws.Range("J362").Value = "Mandarina"
$ws.Range("K362").Value = "Murcott"
$ws.Range("L362").Value = "Primera"
$ws.Range("M362").Value = 110
$ws.Range("N362").Value = 10000
$ws.Range("O362").Value = 10000
$ws.Range("P362").Value = 10000
$ws.Range("Q362").Value = "$/caja 18 kilos"
$ws.Range("R362").Value = "Región de O'Higgins"
$ws.Range("S362").Value = 556
$ws.Range("T362").Value = 18
